$d = $word.ActiveDocument

# --- Edit 1: "Real estate API's ... time)" list item ---
# Change trailing "time)" to "time" and append ", " + "Demographic data "
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("time)")
if ($found) {
    $start = $rng.Start
    $end = $rng.End
    # Remove the trailing ")" character, leaving "time"
    $parenRange = $d.Range($start + 4, $end)
    $parenRange.Text = ""
    # Insert the two new runs right after "time"
    $insPt1 = $d.Range($start + 4, $start + 4)
    $insPt1.InsertAfter(", ")
    $insPt2 = $d.Range($start + 6, $start + 6)
    $insPt2.InsertAfter("Demographic data ")
}

# --- Edit 2: remove the standalone "Demographic data" bullet paragraph   ---
# and make the following "Task Breakdown:" paragraph bold                ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Demographic data") {
        $p.Range.Delete()
        break
    }
}

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Task Breakdown:") {
        $p.Range.Font.Bold = 1
        $p.Range.Font.BoldBi = 1
        break
    }
}

# --- Edit 3: add "Scrum Master: Emmanuel George" to the third empty paragraph ---
# right after the table. That paragraph is the one immediately before the very
# last (struck-through) paragraph of the document.                          ---
$total = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($total)
$targetIndex = $total - 1
if ($lastPara.Range.Font.StrikeThrough -ne -1) {
    # Fallback: search backward for the struck-through paragraph explicitly
    for ($k = $total; $k -ge 1; $k--) {
        if ($d.Paragraphs($k).Range.Font.StrikeThrough -eq -1) {
            $targetIndex = $k - 1
            break
        }
    }
}
$target = $d.Paragraphs($targetIndex)
$target.Range.InsertAfter("S")
$target2 = $d.Paragraphs($targetIndex)
$ins = $d.Range($target2.Range.End - 1, $target2.Range.End - 1)
$ins.InsertAfter("crum Master: Emmanuel George")
